$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 643 (shifts existing rows 643.. down by one),
# matching the diff which adds a new data point (2026/01/13, 火, 20, 201)
# right after the existing 2026/01/13 entries and before the 2026/12/29 block.
$ws.Rows("643:643").Insert()

$dateCell = $ws.Cells.Item(643, 1)
# Force the date column to stay a literal text string (as in the rest of the
# sheet) instead of being auto-parsed into a date serial number.
$dateCell.NumberFormat = "@"
$dateCell.Value = "2026/01/13"
$dateCell.Style = "Normal"

$ws.Cells.Item(643, 2).Value = "火"
$ws.Cells.Item(643, 3).Value = 20
$ws.Cells.Item(643, 4).Value = 201
